$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1900
$ws1.Range("F13").Value = 4494
$ws1.Range("F21").Value = 2310
$ws1.Range("F26").Value = 2204
$ws1.Range("F32").Value = 40

# Sheet "全部类型" (sheet4) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1900
$ws4.Range("F14").Value = 4494
$ws4.Range("F22").Value = 2310
$ws4.Range("F27").Value = 2204
$ws4.Range("F33").Value = 40
